$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.181.37"
$ws.Cells.Item(2, 5).Value = "  -0.83%  "

$ws.Cells.Item(3, 4).Value = "'1.861.00"

$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

$ws.Cells.Item(5, 4).Value = "'0.7080"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "

$ws.Cells.Item(6, 4).Value = "'241.08"
$ws.Cells.Item(6, 5).Value = "  -0.24%  "

$ws.Cells.Item(7, 4).Value = "'1.001"
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

$ws.Cells.Item(8, 5).Value = "  -0.86%  "

$ws.Cells.Item(9, 4).Value = "'0.07645"
$ws.Cells.Item(9, 5).Value = "  -2.37%  "

$ws.Cells.Item(10, 4).Value = "'24.70"
$ws.Cells.Item(10, 5).Value = "  -1.72%  "

$ws.Cells.Item(11, 4).Value = "'0.08423"
$ws.Cells.Item(11, 5).Value = "  +2.28%  "

$ws.Cells.Item(12, 4).Value = "'1.871.33"
$ws.Cells.Item(12, 5).Value = "  -0.33%  "

$ws.Cells.Item(13, 4).Value = "'5.184"
$ws.Cells.Item(13, 5).Value = "  -1.37%  "

$ws.Cells.Item(14, 5).Value = "  -2.42%  "

$ws.Cells.Item(15, 4).Value = "'91.22"
$ws.Cells.Item(15, 5).Value = "  +0.49%  "

$ws.Cells.Item(16, 4).Value = "'29.225.44"
$ws.Cells.Item(16, 5).Value = "  -0.69%  "

$ws.Cells.Item(17, 4).Value = "'5.928"
$ws.Cells.Item(17, 5).Value = "  +0.45%  "

$ws.Cells.Item(18, 4).Value = "'242.51"
$ws.Cells.Item(18, 5).Value = "  -1.94%  "

$ws.Cells.Item(19, 4).Value = "'0.000007809"
$ws.Cells.Item(19, 5).Value = "  -0.71%  "

$ws.Cells.Item(20, 4).Value = "'2.114.59"
$ws.Cells.Item(20, 5).Value = "  +0.07%  "

$ws.Cells.Item(21, 5).Value = "  -1.33%  "

$ws.Cells.Item(22, 4).Value = "'1.000"
$ws.Cells.Item(22, 5).Value = "  +0.07%  "

$ws.Cells.Item(23, 4).Value = "'7.880"
$ws.Cells.Item(23, 5).Value = "  -1.22%  "

$ws.Cells.Item(24, 4).Value = "'1.000"
$ws.Cells.Item(24, 5).Value = "  +0.06%  "

$ws.Cells.Item(26, 4).Value = "'163.21"
$ws.Cells.Item(26, 5).Value = "  -0.35%  "

$ws.Cells.Item(27, 4).Value = "'8.916"
$ws.Cells.Item(27, 5).Value = "  -0.88%  "

$ws.Cells.Item(28, 5).Value = "  +0.85%  "

$ws.Cells.Item(29, 4).Value = "'1.499"
$ws.Cells.Item(29, 5).Value = "  +0.27%  "

$ws.Cells.Item(30, 4).Value = "'1.313"
$ws.Cells.Item(30, 5).Value = "  -3.74%  "

$ws.Cells.Item(31, 5).Value = "  +0.89%  "

$ws.Cells.Item(32, 4).Value = "'4.216"
$ws.Cells.Item(32, 5).Value = "  +2.22%  "

$ws.Cells.Item(33, 4).Value = "'0.05119"
$ws.Cells.Item(33, 5).Value = "  -3.53%  "

$ws.Cells.Item(34, 4).Value = "'0.8064"
$ws.Cells.Item(34, 5).Value = "  +11.69%  "

$ws.Cells.Item(35, 4).Value = "'1.903"
$ws.Cells.Item(35, 5).Value = "  -1.26%  "

$ws.Cells.Item(36, 4).Value = "'1.162"
$ws.Cells.Item(36, 5).Value = "  -3.07%  "

$ws.Cells.Item(37, 4).Value = "'2.681"
$ws.Cells.Item(37, 5).Value = "  +0.12%  "

$ws.Cells.Item(38, 5).Value = "  -0.97%  "

$ws.Cells.Item(39, 4).Value = "'2.699"
$ws.Cells.Item(39, 5).Value = "  -1.02%  "

$ws.Cells.Item(40, 4).Value = "'1.166.33"
$ws.Cells.Item(40, 5).Value = "  -6.30%  "

$ws.Cells.Item(41, 4).Value = "'6.181"
$ws.Cells.Item(41, 5).Value = "  +0.57%  "

$ws.Cells.Item(42, 4).Value = "'0.8917"
$ws.Cells.Item(42, 5).Value = "  -1.77%  "

$ws.Cells.Item(43, 4).Value = "'72.68"
$ws.Cells.Item(43, 5).Value = "  -1.55%  "

$ws.Cells.Item(44, 4).Value = "'1.000"
$ws.Cells.Item(44, 5).Value = "  -0.05%  "

$ws.Cells.Item(45, 4).Value = "'102.02"
$ws.Cells.Item(45, 5).Value = "  -1.12%  "

$ws.Cells.Item(46, 4).Value = "'2.012.05"
$ws.Cells.Item(46, 5).Value = "  -0.05%  "

$ws.Cells.Item(47, 4).Value = "'0.5177"
$ws.Cells.Item(47, 5).Value = "  -2.75%  "

$ws.Cells.Item(48, 4).Value = "'1.768"
$ws.Cells.Item(48, 5).Value = "  -0.11%  "

$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).Value = "'9.249"
$ws.Cells.Item(49, 5).Value = "  -0.05%  "

$ws.Cells.Item(50, 2).Value = "Frax"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(50, 4).Value = "'1.001"
$ws.Cells.Item(50, 5).Value = "  +0.15%  "

$ws.Cells.Item(51, 2).Value = "TheSandbox"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(51, 4).Value = "'0.4266"
$ws.Cells.Item(51, 5).Value = "  -1.08%  "
